# Update the header labels (row 1, columns E:N) on every worksheet of the
# workbook to spell out "severity level(s)" instead of the abbreviated form.
#
# Old -> New
# % 1-2              -> % severity levels 1-2
# # 1-2              -> # severity levels 1-2
# % 3                -> % severity level 3
# # 3                -> # severity level 3
# % 4                -> % severity level 4
# # 4                -> # severity level 4
# % 5                -> % severity level 5
# # 5                -> # severity level 5
# % Tot PiN (3+)     -> % Tot PiN (severity levels 3-5)
# # Tot PiN (3+)     -> # Tot PiN (severity levels 3-5)

$wb = $excel.ActiveWorkbook

$headerMap = @{
    "% 1-2"          = "% severity levels 1-2"
    "# 1-2"          = "# severity levels 1-2"
    "% 3"            = "% severity level 3"
    "# 3"            = "# severity level 3"
    "% 4"            = "% severity level 4"
    "# 4"            = "# severity level 4"
    "% 5"            = "% severity level 5"
    "# 5"            = "# severity level 5"
    "% Tot PiN (3+)" = "% Tot PiN (severity levels 3-5)"
    "# Tot PiN (3+)" = "# Tot PiN (severity levels 3-5)"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($col in 5..14) {
        $cell = $ws.Cells.Item(1, $col)
        $current = $cell.Value2
        if ($headerMap.ContainsKey($current)) {
            $cell.Value = $headerMap[$current]
        }
    }
}
